$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3 with new cluster counts
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 3

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 3

# Add new row 4 - copy formatting (border/bold/center) from A3 since it
# belongs to the same "cluster id" column as A2:A3
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1
